# Actualización automática 2025-11-13 08:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO -------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M2").Value = 1097.67
$wsGrupo.Range("L11").Value = 3036.24
$wsGrupo.Range("M11").Value = 7195.12
$wsGrupo.Range("D30").Value = 457.92

# Row 60 holds "<n> de 58" completion counters per product column.
$wsGrupo.Range("D60").Value = "1 de 58"
$wsGrupo.Range("L60").Value = "1 de 58"
$wsGrupo.Range("M60").Value = "4 de 58"

# --- Sheet: VENTA MENSUAL ----------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F2").Value = 1097.67
$wsMensual.Range("F11").Value = 12403.46
$wsMensual.Range("F30").Value = 457.92
$wsMensual.Range("F60").Value = 19494.86

# --- Sheet: CUMPLIMIENTO MENSUAL ---------------------------------------------
$wsCumplim = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column D (VENTA) widened from 13 to 14 OOXML width units (~13.17 chars).
$wsCumplim.Columns.Item(4).ColumnWidth = 13.166666666666666

$wsCumplim.Range("D3").Value = 457.92
$wsCumplim.Range("E3").Value = 1409.77
$wsCumplim.Range("F3").Value = 0.2451798746044579

$wsCumplim.Range("D11").Value = 3036.24
$wsCumplim.Range("E11").Value = -1595.32
$wsCumplim.Range("F11").Value = 2.107153762873719

$wsCumplim.Range("D12").Value = 13799.11
$wsCumplim.Range("E12").Value = 34241.89
$wsCumplim.Range("F12").Value = 0.2872361108220062

$wsCumplim.Range("D14").Value = 19494.86
$wsCumplim.Range("E14").Value = 38392.49196497848
$wsCumplim.Range("F14").Value = 0.3367723576610359
